$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have a near-empty title row (row 1) above the real
# header row (row 2) and the single data row (row 3). Drop that title
# row so the header moves to row 1 and the data moves to row 2, then
# clear the leftover row-level formatting so we can rebuild it cleanly.
$ws.Rows("1").Delete()
$ws.Rows("1:2").ClearFormats()

# Shift number changed from 2 to 3 (also reflected in the sheet name),
# and the Engineer column now also lists "Usuario de BC".
$ws.Range("B2").Value = 3
$ws.Range("E2").Value = "Renato Hacel Cal y Mayor Rodríguez, Usuario de BC"

# Widen the Activities / Description / Engineer columns. (ColumnWidth is
# expressed in Calibri-11 "characters"; Excel stores width = ColumnWidth +
# 5/6 in the sheet XML, so back that offset out to land on round numbers.)
$ws.Columns("C").ColumnWidth = 29.166666666666668
$ws.Columns("D").ColumnWidth = 49.166666666666664
$ws.Columns("E").ColumnWidth = 29.166666666666668

# Header row: bold white text on the existing blue fill, thin border,
# centered both ways.
$hdr = $ws.Range("A1:E1")
$hdr.Font.Bold = $true
$hdr.Font.Color = 16777215
$hdr.Interior.Color = 12611584
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108

# Data row: thin border, wrapped text anchored to the top.
$body = $ws.Range("A2:E2")
$body.Borders.LineStyle = 1
$body.VerticalAlignment = -4160
$body.WrapText = $true

# Reflect the new shift number in the sheet name/title.
$ws.Name = "Shift 3 - 2025-02-16"
